$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: update market data file name + param dict (C1, D1) ---
$ws.Range("C1").Value = "data/market_data2.csv"
$ws.Range("D1").Value = "{'LOT_SIZE': 11, 'HIST_LENGTH': 10, 'delta': 1, 'preferred_lots_low': 8, 'preferred_lots_high': 12, 'margin_adjustment': 0.0005, 'MarketDataFile': 'data/market_data2.csv'}"

# --- Row 3: updated benchmark numbers for humming_trader ---
$ws.Range("B3").Value = 4257
$ws.Range("C3").Value = 4354
$ws.Range("D3").Value = -97
$ws.Range("E3").Value = 97
$ws.Range("H3").Value = -62791
$ws.Range("I3").Value = 1241491
$ws.Range("J3").Value = 1251191

# --- Remove the old "Match #2" block (rows 6-8) ---
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()

# --- Add new row 4 with optiver_trader results ---
$ws.Range("A4").Value = "optiver_trader"
$ws.Range("B4").Value = 6280
$ws.Range("C4").Value = 6284
$ws.Range("D4").Value = -4
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 134200
$ws.Range("G4").Value = 134300
$ws.Range("H4").Value = 94226
$ws.Range("I4").Value = -549426
$ws.Range("J4").Value = -549026
$ws.Range("K4").Value = "OK"

# New fill/style (red) applied to the ProfitOrLoss cell of the new row
$ws.Range("J4").Interior.Color = 255

# --- Fix the used range/dimension ---
$ws.UsedRange | Out-Null
